{"js": "// \"Remove revision and date\"\n//\n// 1. Delete the \"Date: ...\" and \"Revision: ...\" paragraphs that sit\n//    between the title and the \"General Formatting Rules\" heading.\n// 2. The heading itself was originally split across two runs (\"General\"\n//    and \" Formatting Rules\", with a bookmark in between); fold them into\n//    a single run reading \"General Formatting Rules\" while leaving the\n//    bookmark in place.\n\nconst body = context.document.body;\n\n// --- Step 1: remove the Date/Revision paragraphs -------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst paragraphsToRemove = [\"Date: 03/10/2016\", \"Revision: Draft\"];\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  const text = paragraphs.items[i].text.trim();\n  if (paragraphsToRemove.indexOf(text) !== -1) {\n    paragraphs.items[i].delete();\n  }\n}\nawait context.sync();\n\n// --- Step 2: merge \"General\" + \" Formatting Rules\" into one run ----------\nconst headingSearch = body.search(\"General Formatting Rules\", { matchCase: true });\nheadingSearch.load(\"text\");\nawait context.sync();\n\nif (headingSearch.items.length > 0) {\n  const headingRange = headingSearch.items[0];\n\n  // Split the matched range on spaces so we can isolate the \"General\"\n  // word (first run) from the rest, without disturbing the bookmark that\n  // sits between the two original runs.\n  const words = headingRange.getTextRanges([\" \"], true);\n  words.load(\"text\");\n  await context.sync();\n\n  // Remove the leading \"General\" run entirely.\n  words.items[0].insertText(\"\", \"Replace\");\n  await context.sync();\n\n  // The remaining run still reads \" Formatting Rules\" (leading space from\n  // the original run); replace it with the fully merged heading text.\n  const remainder = body.search(\" Formatting Rules\", { matchCase: true });\n  remainder.load(\"text\");\n  await context.sync();\n  if (remainder.items.length > 0) {\n    remainder.items[0].insertText(\"General Formatting Rules\", \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# \"Remove revision and date\"\n#\n# 1. Delete the \"Date: ...\" and \"Revision: ...\" paragraphs that sit\n#    between the title and the \"General Formatting Rules\" heading.\n# 2. The heading itself was originally split across two runs (\"General\"\n#    and \" Formatting Rules\", with a bookmark in between); fold them into\n#    a single run reading \"General Formatting Rules\" while leaving the\n#    bookmark in place.\n\n$d = $word.ActiveDocument\n\nfunction Remove-ParagraphByText($doc, $text) {\n    $range = $doc.Content\n    $range.Find.ClearFormatting()\n    $found = $range.Find.Execute($text)\n    if ($found) {\n        $range.Expand(4) | Out-Null   # wdParagraph = 4: grow to the whole paragraph (incl. mark)\n        $range.Delete()\n    }\n    return $found\n}\n\n# --- Step 1: remove the Date/Revision paragraphs --------------------------\nRemove-ParagraphByText $d \"Date: 03/10/2016\" | Out-Null\nRemove-ParagraphByText $d \"Revision: Draft\" | Out-Null\n\n# --- Step 2: merge \"General\" + \" Formatting Rules\" into one run -----------\n# Locate the heading paragraph (Word's Range.Text concatenates run text, so\n# the split still reads as one phrase here).\n$headingParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*General*Formatting Rules*\") {\n        $headingParagraph = $p\n        break\n    }\n}\n\nif ($headingParagraph -ne $null) {\n    # Only touch it if \"General\" is still its own leading run (i.e. there is\n    # something right after it, such as the bookmark, before the rest of the\n    # text) -- keeps this step a no-op if re-applied to an already-merged\n    # heading instead of corrupting the text.\n    $leadRange = $headingParagraph.Range.Duplicate()\n    $leadRange.Find.ClearFormatting()\n    $leadRange.Find.Execute(\"General\") | Out-Null\n\n    $afterGeneral = $headingParagraph.Range.Duplicate()\n    $afterGeneral.Start = $leadRange.End\n    $afterGeneral.End = $headingParagraph.Range.End\n\n    if ($afterGeneral.Text -like \" Formatting Rules*\") {\n        # Remove just the leading \"General\" run/word, leaving the bookmark\n        # that sits right after it untouched.\n        $leadRange.Text = \"\"\n\n        # The remaining run now reads \" Formatting Rules\" (leading space\n        # left over from the original run); replace it with the merged text.\n        $fixRange = $headingParagraph.Range.Duplicate()\n        $fixRange.Find.ClearFormatting()\n        $fixRange.Find.Execute(\" Formatting Rules\") | Out-Null\n        $fixRange.Text = \"General Formatting Rules\"\n    }\n}\n"}
